{"js": "// Update the date line and every answer cell in the practice table with\n// the new day's values (see commit: \"Update master to output generated\n// at c8c62b6\").\n\n// --- 1) Update the date paragraph (first paragraph in the body) -------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2025-07-13 Sunday\", Word.InsertLocation.replace);\n\n// --- 2) Update every cell of the 20x5 answer table --------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"7+24=31\", \"96-9=87\", \"95-37=58\", \"91-63=28\", \"94-88=6\"],\n  [\"25+48=73\", \"68-49=19\", \"80-17=63\", \"9+13=22\", \"17+55=72\"],\n  [\"33-18=15\", \"25+57=82\", \"72-49=23\", \"12-7=5\", \"15+19=34\"],\n  [\"93-79=14\", \"90-24=66\", \"60-8=52\", \"94-59=35\", \"85-76=9\"],\n  [\"91-9=82\", \"62+29=91\", \"86-47=39\", \"92-79=13\", \"19+8=27\"],\n  [\"39+44=83\", \"87+9=96\", \"66+9=75\", \"37+44=81\", \"38+9=47\"],\n  [\"54-15=39\", \"9+13=22\", \"40-12=28\", \"6+78=84\", \"44+9=53\"],\n  [\"84-29=55\", \"81-25=56\", \"90-24=66\", \"60-35=25\", \"29+45=74\"],\n  [\"38+24=62\", \"4+49=53\", \"22+49=71\", \"90-18=72\", \"96-37=59\"],\n  [\"69+25=94\", \"17+9=26\", \"28+39=67\", \"72-65=7\", \"93-26=67\"],\n  [\"4+79=83\", \"33-29=4\", \"67-38=29\", \"74-39=35\", \"90-8=82\"],\n  [\"69+13=82\", \"72-29=43\", \"94-25=69\", \"87-78=9\", \"18+43=61\"],\n  [\"84-75=9\", \"73-29=44\", \"51-14=37\", \"86-49=37\", \"61-19=42\"],\n  [\"35+49=84\", \"75-26=49\", \"77+18=95\", \"35-19=16\", \"70-36=34\"],\n  [\"19+27=46\", \"77-68=9\", \"51-28=23\", \"73-8=65\", \"7+89=96\"],\n  [\"50-47=3\", \"79+12=91\", \"27+69=96\", \"41-29=12\", \"8+77=85\"],\n  [\"48+24=72\", \"18+26=44\", \"16+9=25\", \"82-6=76\", \"47+16=63\"],\n  [\"13+78=91\", \"61-14=47\", \"19+68=87\", \"37+18=55\", \"23+9=32\"],\n  [\"6+77=83\", \"92-16=76\", \"16+15=31\", \"51-7=44\", \"71-53=18\"],\n  [\"52+19=71\", \"65-58=7\", \"25+46=71\", \"43-7=36\", \"56-48=8\"],\n];\n\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Update the date line and every answer cell in the practice table with\n# the new day's values (see commit: \"Update master to output generated\n# at c8c62b6\").\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date paragraph (first paragraph in the body) -------\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateParagraph.Range.Text = \"2025-07-13 Sunday\"\n\n# --- 2) Update every cell of the 20x5 answer table --------------------\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(\"7+24=31\", \"96-9=87\", \"95-37=58\", \"91-63=28\", \"94-88=6\"),\n    @(\"25+48=73\", \"68-49=19\", \"80-17=63\", \"9+13=22\", \"17+55=72\"),\n    @(\"33-18=15\", \"25+57=82\", \"72-49=23\", \"12-7=5\", \"15+19=34\"),\n    @(\"93-79=14\", \"90-24=66\", \"60-8=52\", \"94-59=35\", \"85-76=9\"),\n    @(\"91-9=82\", \"62+29=91\", \"86-47=39\", \"92-79=13\", \"19+8=27\"),\n    @(\"39+44=83\", \"87+9=96\", \"66+9=75\", \"37+44=81\", \"38+9=47\"),\n    @(\"54-15=39\", \"9+13=22\", \"40-12=28\", \"6+78=84\", \"44+9=53\"),\n    @(\"84-29=55\", \"81-25=56\", \"90-24=66\", \"60-35=25\", \"29+45=74\"),\n    @(\"38+24=62\", \"4+49=53\", \"22+49=71\", \"90-18=72\", \"96-37=59\"),\n    @(\"69+25=94\", \"17+9=26\", \"28+39=67\", \"72-65=7\", \"93-26=67\"),\n    @(\"4+79=83\", \"33-29=4\", \"67-38=29\", \"74-39=35\", \"90-8=82\"),\n    @(\"69+13=82\", \"72-29=43\", \"94-25=69\", \"87-78=9\", \"18+43=61\"),\n    @(\"84-75=9\", \"73-29=44\", \"51-14=37\", \"86-49=37\", \"61-19=42\"),\n    @(\"35+49=84\", \"75-26=49\", \"77+18=95\", \"35-19=16\", \"70-36=34\"),\n    @(\"19+27=46\", \"77-68=9\", \"51-28=23\", \"73-8=65\", \"7+89=96\"),\n    @(\"50-47=3\", \"79+12=91\", \"27+69=96\", \"41-29=12\", \"8+77=85\"),\n    @(\"48+24=72\", \"18+26=44\", \"16+9=25\", \"82-6=76\", \"47+16=63\"),\n    @(\"13+78=91\", \"61-14=47\", \"19+68=87\", \"37+18=55\", \"23+9=32\"),\n    @(\"6+77=83\", \"92-16=76\", \"16+15=31\", \"51-7=44\", \"71-53=18\"),\n    @(\"52+19=71\", \"65-58=7\", \"25+46=71\", \"43-7=36\", \"56-48=8\")\n)\n\nfor ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n\nWrite-Output \"done\"\n"}
